{"js": "// Find the \"Requisitos\" answer paragraph (\"LOT2007: Bioqu\u00edmica I (Indica\u00e7\u00e3o de Conjunto)\")\n// and remove the three paragraphs that followed it in the old footer block:\n//   1. a blank paragraph\n//   2. \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3. \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n// The paragraph with \"LOT2007...\" itself, and the blank paragraph that follows\n// the removed block (right before the trailing page-break paragraph), are kept.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst marker = \"LOT2007: Bioqu\u00edmica I (Indica\u00e7\u00e3o de Conjunto)\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === marker) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error('Could not find paragraph with text \"' + marker + '\"');\n}\n\n// Collect the next three paragraphs (blank, \"Ver no Jupiter...\", \"\u00a9 2020...\")\nconst toDelete = [];\nlet cursor = anchor;\nfor (let i = 0; i < 3; i++) {\n  cursor = cursor.getNext();\n  toDelete.push(cursor);\n}\ntoDelete.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\n// Sanity-check the texts before deleting, then delete them.\nconst expected = [\n  \"\",\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\",\n];\nfor (let i = 0; i < toDelete.length; i++) {\n  if (toDelete[i].text !== expected[i]) {\n    throw new Error(\n      \"Unexpected paragraph text at offset \" + (i + 1) + \": \" + JSON.stringify(toDelete[i].text)\n    );\n  }\n}\n\ntoDelete.forEach((p) => p.delete());\nawait context.sync();\n", "ps1": "# Find the \"Requisitos\" answer paragraph (\"LOT2007: Bioqu\u00edmica I (Indica\u00e7\u00e3o de Conjunto)\")\n# and remove the three paragraphs that followed it in the old footer block:\n#   1. a blank paragraph\n#   2. \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3. \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n# The paragraph with \"LOT2007...\" itself, and the blank paragraph that follows\n# the removed block (right before the trailing page-break paragraph), are kept.\n\n$d = $word.ActiveDocument\n\n$marker = \"LOT2007: Bioqu\u00edmica I (Indica\u00e7\u00e3o de Conjunto)\"\n$anchorIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $txt = $d.Paragraphs.Item($i).Range.Text\n    $txt = $txt.TrimEnd(\"`r\", \"`a\")\n    if ($txt -eq $marker) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find paragraph with text '$marker'\"\n}\n\n$expected = @(\n    \"\",\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n)\n\n# Deleting always targets the paragraph right after the anchor, since earlier\n# deletions shift the following ones up into that slot.\nfor ($k = 0; $k -lt 3; $k++) {\n    $p = $d.Paragraphs.Item($anchorIndex + 1)\n    $txt = $p.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($txt -ne $expected[$k]) {\n        throw \"Unexpected paragraph text at offset $($k + 1): '$txt'\"\n    }\n    $p.Range.Delete()\n}\n"}
